$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The "SDA1 (AIC #2)" textbox ("TextBox 8") lives nested inside the
# top-level group "Group 41" (GroupItems flattens the inner "Group 39").
$grp = $s.Shapes.Item(1)
$tb = $grp.GroupItems.Item(7)

$tr = $tb.TextFrame.TextRange
# Split off the leading "SDA1 " (5 chars, including the trailing space)
# and replace it with "SDA2 ", leaving "(AIC #2)" as a separate run with
# the same formatting.
$lead = $tr.Characters(1, 5)
$lead.Text = "SDA2 "
